$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Arveja Verde" (Vega Modelo de
# Temuco). Insert a fresh row at 12 — this pushes the existing rows 12-76
# down to 13-77 (and the sheet's used range grows from R76 to R77) — then
# fill it in with the new observation's data.
$ws.Rows.Item(12).EntireRow.Insert()

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(12, 3).Value = "La Araucanía"
$ws.Cells.Item(12, 4).Value = 44550
$ws.Cells.Item(12, 5).Value = 9
$ws.Cells.Item(12, 6).Value = 100112022
$ws.Cells.Item(12, 7).Value = "Arveja Verde"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 95
$ws.Cells.Item(12, 11).Value = 15000
$ws.Cells.Item(12, 12).Value = 15000
$ws.Cells.Item(12, 13).Value = 15000
$ws.Cells.Item(12, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(12, 16).Value = 600
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"
